$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '245.93'
$ws.Range('D2').Style = 'Normal'

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '22.19'
$ws.Range('D3').Style = 'Normal'

$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '5.346'
$ws.Range('D4').Style = 'Normal'

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '0.05922'
$ws.Range('D5').Style = 'Normal'

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '6.391'
$ws.Range('D7').Style = 'Normal'

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.8128'
$ws.Range('D8').Style = 'Normal'

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.9607'
$ws.Range('D9').Style = 'Normal'

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.1428'
$ws.Range('D10').Style = 'Normal'

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07404'
$ws.Range('D11').Style = 'Normal'

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.03429'
$ws.Range('D12').Style = 'Normal'

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.03041'
$ws.Range('D13').Style = 'Normal'

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '4.006'
$ws.Range('D14').Style = 'Normal'

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.09390'
$ws.Range('D15').Style = 'Normal'

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.001588'
$ws.Range('D16').Style = 'Normal'

$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '17OneONEWorstin24h'
$ws.Range('E18').Style = 'Normal'

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.006178'
$ws.Range('D19').Style = 'Normal'

$ws.Range('B20').NumberFormat = '@'
$ws.Range('B20').Value = 'UpBots'
$ws.Range('B20').Style = 'Normal'

$ws.Range('C20').NumberFormat = '@'
$ws.Range('C20').Value = 'https://coinranking.com/coin/m5ozaAIK6+upbots-ubxt'
$ws.Range('C20').Style = 'Normal'

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.007490'
$ws.Range('D20').Style = 'Normal'

$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '19UpBotsUBXTBestin24h'
$ws.Range('E20').Style = 'Normal'

$ws.Range('B21').NumberFormat = '@'
$ws.Range('B21').Value = 'HotbitToken'
$ws.Range('B21').Style = 'Normal'

$ws.Range('C21').NumberFormat = '@'
$ws.Range('C21').Value = 'https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb'
$ws.Range('C21').Style = 'Normal'

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.004085'
$ws.Range('D21').Style = 'Normal'

$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '20HotbitTokenHTB'
$ws.Range('E21').Style = 'Normal'

$ws.Range('B22').NumberFormat = '@'
$ws.Range('B22').Value = 'BitKan'
$ws.Range('B22').Style = 'Normal'

$ws.Range('C22').NumberFormat = '@'
$ws.Range('C22').Value = 'https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan'
$ws.Range('C22').Style = 'Normal'

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.0009843'
$ws.Range('D22').Style = 'Normal'

$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '21BitKanKAN'
$ws.Range('E22').Style = 'Normal'

$ws.Range('B23').NumberFormat = '@'
$ws.Range('B23').Value = 'NitroEx'
$ws.Range('B23').Style = 'Normal'

$ws.Range('C23').NumberFormat = '@'
$ws.Range('C23').Value = 'https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx'
$ws.Range('C23').Style = 'Normal'

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.00009703'
$ws.Range('D23').Style = 'Normal'

$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '22NitroExNTX'
$ws.Range('E23').Style = 'Normal'

$ws.Range('B24').NumberFormat = '@'
$ws.Range('B24').Value = 'LEO'
$ws.Range('B24').Style = 'Normal'

$ws.Range('C24').NumberFormat = '@'
$ws.Range('C24').Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range('C24').Style = 'Normal'

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '3.741'
$ws.Range('D24').Style = 'Normal'

$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '23LEOLEO'
$ws.Range('E24').Style = 'Normal'

$ws.Range('B25').NumberFormat = '@'
$ws.Range('B25').Value = 'BTSEToken'
$ws.Range('B25').Style = 'Normal'

$ws.Range('C25').NumberFormat = '@'
$ws.Range('C25').Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range('C25').Style = 'Normal'

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.165'
$ws.Range('D25').Style = 'Normal'

$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '24BTSETokenBTSE'
$ws.Range('E25').Style = 'Normal'

$ws.Range('B26').NumberFormat = '@'
$ws.Range('B26').Value = 'BitpandaEcosystemToken'
$ws.Range('B26').Style = 'Normal'

$ws.Range('C26').NumberFormat = '@'
$ws.Range('C26').Value = 'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best'
$ws.Range('C26').Style = 'Normal'

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.3268'
$ws.Range('D26').Style = 'Normal'

$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '25BitpandaEcosystemTokenBEST'
$ws.Range('E26').Style = 'Normal'

$ws.Range('B27').NumberFormat = '@'
$ws.Range('B27').Value = 'ProBitToken'
$ws.Range('B27').Style = 'Normal'

$ws.Range('C27').NumberFormat = '@'
$ws.Range('C27').Value = 'https://coinranking.com/coin/lQP4d6T2+probittoken-prob'
$ws.Range('C27').Style = 'Normal'

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.1333'
$ws.Range('D27').Style = 'Normal'

$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '26ProBitTokenPROB'
$ws.Range('E27').Style = 'Normal'

$ws.Range('B41').NumberFormat = '@'
$ws.Range('B41').Value = 'KickToken'
$ws.Range('B41').Style = 'Normal'

$ws.Range('C41').NumberFormat = '@'
$ws.Range('C41').Value = 'https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick'
$ws.Range('C41').Style = 'Normal'

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.006472'
$ws.Range('D41').Style = 'Normal'

$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '40KickTokenKICK'
$ws.Range('E41').Style = 'Normal'

$ws.Range('B42').NumberFormat = '@'
$ws.Range('B42').Value = 'BKEXToken'
$ws.Range('B42').Style = 'Normal'

$ws.Range('C42').NumberFormat = '@'
$ws.Range('C42').Value = 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'
$ws.Range('C42').Style = 'Normal'

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.1072'
$ws.Range('D42').Style = 'Normal'

$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '41BKEXTokenBKK'
$ws.Range('E42').Style = 'Normal'

$ws.Range('B43').NumberFormat = '@'
$ws.Range('B43').Value = 'CEJI'
$ws.Range('B43').Style = 'Normal'

$ws.Range('C43').NumberFormat = '@'
$ws.Range('C43').Value = 'https://coinranking.com/coin/SbKjCVJCh+ceji-ceji'
$ws.Range('C43').Style = 'Normal'

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.003001'
$ws.Range('D43').Style = 'Normal'

$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '42CEJICEJI'
$ws.Range('E43').Style = 'Normal'

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.005921'
$ws.Range('D44').Style = 'Normal'

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.00005297'
$ws.Range('D45').Style = 'Normal'

$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '46CoinbaseStockTokenCOIN'
$ws.Range('E47').Style = 'Normal'

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.04728'
$ws.Range('D48').Style = 'Normal'

$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '47BOLOBOLO'
$ws.Range('E48').Style = 'Normal'

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.00002101'
$ws.Range('D49').Style = 'Normal'
